$wb = $excel.ActiveWorkbook

# --- Trucks sheet: capacity values 100000 -> 250 (D2:D11) ---
$wsTrucks = $wb.Worksheets.Item("Trucks")
$wsTrucks.Range("D2:D11").Value = 250

# --- Selection / active-cell bookkeeping on the sheets touched by the edit ---

# Customers: H6:H7 (H7 active) -> D2:D7 (D2 active), no longer the tab shown
$wsCustomers = $wb.Worksheets.Item("Customers")
$wsCustomers.Activate() | Out-Null
$wsCustomers.Range("D2:D7").Select() | Out-Null

# Nodes: H18 -> D17
$wsNodes = $wb.Worksheets.Item("Nodes")
$wsNodes.Activate() | Out-Null
$wsNodes.Range("D17").Select() | Out-Null

# Others: B6 -> E1:I8 (E1 active)
$wsOthers = $wb.Worksheets.Item("Others")
$wsOthers.Activate() | Out-Null
$wsOthers.Range("E1:I8").Select() | Out-Null

# --- New sheet "Tabelle1" appended after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNew = $wb.Worksheets.Add($null, $lastSheet)
$wsNew.Name = "Tabelle1"
$wsNew.Range("A1:E9").Select() | Out-Null

# Trucks: B3:B11 (B3 active) -> D3:D11 (D3 active); becomes the selected tab (activeTab=2).
# Activating it last ensures it (not the freshly-added Tabelle1) ends up as the
# workbook's tabSelected / activeTab sheet, matching the target state.
$wsTrucks.Activate() | Out-Null
$wsTrucks.Range("D3:D11").Select() | Out-Null
